# "added simple plot option" — refresh the sample data table (A2:D9) with a
# new data set (8 rows instead of 6) and tidy up the view/print state the
# way Excel would after someone reworked the sheet to drive a quick chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data table: columns are y, x, y_err, x_err (headers already in row 1) ---

# Row 2
$ws.Cells.Item(2, 1).Value = -0.777
$ws.Cells.Item(2, 2).Value = 450
$ws.Cells.Item(2, 3).Value = 0.001
$ws.Cells.Item(2, 4).Value = 0

# Row 3
$ws.Cells.Item(3, 1).Value = -0.551
$ws.Cells.Item(3, 2).Value = 500
$ws.Cells.Item(3, 3).Value = 0.001
$ws.Cells.Item(3, 4).Value = 0.01

# Row 4
$ws.Cells.Item(4, 1).Value = -0.529
$ws.Cells.Item(4, 2).Value = 503
$ws.Cells.Item(4, 3).Value = 0.001
$ws.Cells.Item(4, 4).Value = 0.01

# Row 5
$ws.Cells.Item(5, 1).Value = -0.427
$ws.Cells.Item(5, 2).Value = 525
$ws.Cells.Item(5, 3).Value = 0.001
$ws.Cells.Item(5, 4).Value = 0.01

# Row 6
$ws.Cells.Item(6, 1).Value = -0.352
$ws.Cells.Item(6, 2).Value = 540
$ws.Cells.Item(6, 3).Value = 0.001
$ws.Cells.Item(6, 4).Value = 0.01

# Row 7
$ws.Cells.Item(7, 1).Value = 0.342
$ws.Cells.Item(7, 2).Value = 552
$ws.Cells.Item(7, 3).Value = 0.001
$ws.Cells.Item(7, 4).Value = 0.01

# Row 8 (new)
$ws.Cells.Item(8, 1).Value = -0.262
$ws.Cells.Item(8, 2).Value = 575
$ws.Cells.Item(8, 3).Value = 0.001
$ws.Cells.Item(8, 4).Value = 0.001

# Row 9 (new)
$ws.Cells.Item(9, 1).Value = -0.1727
$ws.Cells.Item(9, 2).Value = 602
$ws.Cells.Item(9, 3).Value = 0.0001
$ws.Cells.Item(9, 4).Value = 0.001

# --- View state: leave the cursor on B2, like after reviewing the new data ---
[void]$ws.Range("B2").Select()

# --- Print setup: page got configured (portrait) while setting up the plot/printout ---
$ws.PageSetup.Orientation = 1

Write-Output "edit applied"
